$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row additions: P1=14, Q1=15, matching header style of existing header cells ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy formatting (style) from an existing header cell (B1) onto the new header cells
$ws.Range("B1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows 2-25: swap I/K and M/O values, add new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new
    $ws.Cells.Item($r, 17).Value = 2   # Q: new
}
